$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Y-GT -> RDW-CV
$ws.Range("A4").Value = "RDW-CV"
$ws.Range("B4").Value = "γ谷氨酰转肽酶"

# Row 8: T-bil -> RDW-CV
$ws.Range("A8").Value = "RDW-CV"
$ws.Range("B8").Value = "总胆红素"
$ws.Range("D8").Value = "μmol/L"

# Row 9: D-bil -> RDW-CV
$ws.Range("A9").Value = "RDW-CV"
$ws.Range("B9").Value = "直接胆红素"
$ws.Range("D9").Value = "μmol/L"

# Row 10: HBsAg -> RDW-CV
$ws.Range("A10").Value = "RDW-CV"

# Row 13: AST/ALT -> RDW-CV
$ws.Range("A13").Value = "RDW-CV"
